$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: "link zoom" header + zoom link value (as plain text) / hyperlink
$url = "https://us05web.zoom.us/j/83990720254?pwd=2BJciSFGDWGEaP2QzQILctCFu0D4Hv.1"

$ws.Range("D1").Value = "link zoom"
$ws.Hyperlinks.Add($ws.Range("D2"), $url)
$ws.Range("D3").Value = $url

# Column D width
$ws.Columns("D").ColumnWidth = 65

# Bold header row
$ws.Range("A1:D1").Font.Bold = $true

# Page orientation
$ws.PageSetup.Orientation = 1

# Selection lands on D3 after the edits
$ws.Range("D3").Select()
